# Update the "stat output" figures on Sheet3:
#   - A2 (minutes per match) goes from 2 to 1
#   - A3 (number of matches) goes from 30 to 20
# The dependent formulas in F1, F3 and F5 recalc automatically.
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("A2").Value = 1
$ws3.Range("A3").Value = 20

# Sheet3 becomes the active/selected sheet and tab, with A2 selected.
$ws3.Activate()
$ws3.Range("A2").Select()

$wb.Save()
